$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 318, shifting existing rows 318-328 down to 319-329.
$ws.Rows.Item(318).Insert()

# Populate the newly inserted row 318 with the new record.
$ws.Cells.Item(318, 1).Value = 9
$ws.Cells.Item(318, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(318, 3).Value = "Metropolitana"
$ws.Cells.Item(318, 4).Value = 45075
$ws.Cells.Item(318, 4).NumberFormat = $ws.Cells.Item(319, 4).NumberFormat
$ws.Cells.Item(318, 5).Value = 13
$ws.Cells.Item(318, 6).Value = 100112026
$ws.Cells.Item(318, 7).Value = "Haba"
$ws.Cells.Item(318, 8).Value = "Sin especificar"
$ws.Cells.Item(318, 9).Value = "Primera"
$ws.Cells.Item(318, 10).Value = 34
$ws.Cells.Item(318, 11).Value = 22000
$ws.Cells.Item(318, 12).Value = 24000
$ws.Cells.Item(318, 13).Value = 23000
$ws.Cells.Item(318, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(318, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(318, 16).Value = 920
$ws.Cells.Item(318, 17).Value = 25
$ws.Cells.Item(318, 18).Value = "Hortaliza"
